# Roadmap.xlsx update — "Before chaning factorio version"
#
# Roadmap sheet: the squad-building related tasks move from to-do/in-progress
# into done, several new sub-tasks are recorded, and "Hunting cabin" becomes
# the new in-progress item.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Roadmap")

# --- Insert 6 new rows right after row 33 (old rows 34.. shift down to 40..) ---
$ws.Range("A34:C39").EntireRow.Insert()

# Row 33 ("Split select unit window from workshop") flips from "in progress" to "done"
$ws.Range("B33").Value = "done"

# The 6 freshly inserted rows (34-39) plus touching up 35/36 (old Sqad
# templates / Build squad window, now reordered) — build out the final
# block of 9 "done" rows (33-41).
$ws.Range("A34").Value = "Split itemlistbuilder from storage ui"
$ws.Range("B34").Value = "done"

$ws.Range("A35").Value = "Sqad templates"
$ws.Range("B35").Value = "done"

$ws.Range("A36").Value = "Build squad window"
$ws.Range("B36").Value = "done"

$ws.Range("A37").Value = "Squad queue window"
$ws.Range("B37").Value = "done"

$ws.Range("A38").Value = "Squad templates window"
$ws.Range("B38").Value = "done"

$ws.Range("A39").Value = "Squad list window"
$ws.Range("B39").Value = "done"

$ws.Range("A40").Value = "Squad window"
$ws.Range("B40").Value = "done"

$ws.Range("A41").Value = "Try to build squads every x seconds"
$ws.Range("B41").Value = "done"

# "Hunting cabin" (now row 42 after the insert) becomes the new in-progress item
$ws.Range("B42").Value = "in progress"

# --- Fix up the conditional formatting ranges that shifted with the insert ---
# The lone single-cell rule block tracked the old "gap" row (B38) and now
# needs to track the new gap row (B44).
$gapRules = $ws.Range("B38").FormatConditions
for ($i = 1; $i -le $gapRules.Count; $i++) {
    $gapRules.Item($i).ModifyAppliesToRange($ws.Range("B44"))
}

# The main rule block covered everything else (B1:B37, B39:B1048576); keep it
# tracking the rest of the status column around the relocated gap row.
$mainRules = $ws.Range("B1").FormatConditions
for ($i = 1; $i -le $mainRules.Count; $i++) {
    $mainRules.Item($i).ModifyAppliesToRange($ws.Range("B1:B43"))
}

# --- Update the view to match where the author left the cursor ---
$ws.Range("B43").Select()
